# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1184.2142
$ws.Range("I19").Value = 726.4545000000001
$ws.Range("J19").Value = 1480.4117
$ws.Range("K19").Value = 726.4545000000001
$ws.Range("L19").Value = 1480.4117
$ws.Range("M19").Value = -551.4545000000001
$ws.Range("N19").Value = -1830.4117

$ws.Range("H88").Value = 24306.312
$ws.Range("I88").Value = 63955.332
$ws.Range("J88").Value = 8791.478999999999
$ws.Range("K88").Value = 63955.332
$ws.Range("L88").Value = 8791.478999999999
$ws.Range("M88").Value = -63549.332
$ws.Range("N88").Value = -9603.478999999999

$ws.Range("H91").Value = 24306.312
$ws.Range("I91").Value = 63955.332
$ws.Range("J91").Value = 8791.478999999999
$ws.Range("K91").Value = 63955.332
$ws.Range("L91").Value = 8791.478999999999
$ws.Range("M91").Value = -62551.332
$ws.Range("N91").Value = -11599.479

$ws.Range("H132").Value = 5648.08
$ws.Range("I132").Value = 1652.1904
$ws.Range("K132").Value = 4956.5712
$ws.Range("M132").Value = -2426.5712

$ws.Range("H135").Value = 47620084
$ws.Range("I135").Value = 1046.8334
$ws.Range("J135").Value = 111112130
$ws.Range("K135").Value = 9421.500599999999
$ws.Range("L135").Value = 1000009170
$ws.Range("M135").Value = -6886.500599999999
$ws.Range("N135").Value = -1000014240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1867.7826
$ws.Range("I2").Value = 2115.0908
$ws.Range("K2").Value = 2115.0908
$ws.Range("M2").Value = -2002.0908

$ws.Range("H88").Value = 2153.476
$ws.Range("I88").Value = 2410.6667
$ws.Range("J88").Value = 1960.5834
$ws.Range("K88").Value = 2410.6667
$ws.Range("L88").Value = 1960.5834
$ws.Range("M88").Value = -2004.6667
$ws.Range("N88").Value = -2772.5834

$ws.Range("H91").Value = 2153.476
$ws.Range("I91").Value = 2410.6667
$ws.Range("J91").Value = 1960.5834
$ws.Range("K91").Value = 2410.6667
$ws.Range("L91").Value = 1960.5834
$ws.Range("M91").Value = -1006.6667
$ws.Range("N91").Value = -4768.5834

$ws.Range("H97").Value = 1120.625
$ws.Range("I97").Value = 994.1667
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 994.1667
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -498.1667
$ws.Range("N97").Value = -2492

$ws.Range("H110").Value = 1229.0667
$ws.Range("I110").Value = 1042.5714
$ws.Range("J110").Value = 1392.25
$ws.Range("K110").Value = 1042.5714
$ws.Range("L110").Value = 1392.25
$ws.Range("M110").Value = 1002.4286
$ws.Range("N110").Value = -5482.25

$ws.Range("H116").Value = 1867.7826
$ws.Range("I116").Value = 2115.0908
$ws.Range("K116").Value = 2115.0908
$ws.Range("M116").Value = 178.9092000000001

$ws.Range("H132").Value = 1759.375
$ws.Range("I132").Value = 1598.3334
$ws.Range("J132").Value = 1966.4286
$ws.Range("K132").Value = 4795.0002
$ws.Range("L132").Value = 5899.2858
$ws.Range("M132").Value = -2265.0002
$ws.Range("N132").Value = -10959.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1867.7826
$ws.Range("I3").Value = 2115.0908
$ws.Range("K3").Value = 2115.0908
$ws.Range("M3").Value = -2001.0908

$ws.Range("H82").Value = 18814.25
$ws.Range("J82").Value = 23333.334
$ws.Range("L82").Value = 23333.334
$ws.Range("N82").Value = -24099.334

$ws.Range("H85").Value = 18814.25
$ws.Range("J85").Value = 23333.334
$ws.Range("L85").Value = 23333.334
$ws.Range("N85").Value = -25985.334

$ws.Range("H94").Value = 2247.1304
$ws.Range("I94").Value = 1793.1765
$ws.Range("J94").Value = 3533.3333
$ws.Range("K94").Value = 1793.1765
$ws.Range("L94").Value = 3533.3333
$ws.Range("M94").Value = -1342.1765
$ws.Range("N94").Value = -4435.3333

$ws.Range("H105").Value = 2486.6155
$ws.Range("I105").Value = 2440.5881
$ws.Range("K105").Value = 2440.5881
$ws.Range("M105").Value = -693.5880999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2172.3286
$ws.Range("I31").Value = 1088.4615
$ws.Range("J31").Value = 3535.9033
$ws.Range("K31").Value = 1088.4615
$ws.Range("L31").Value = 3535.9033
$ws.Range("M31").Value = -793.4614999999999
$ws.Range("N31").Value = -4125.9033

$ws.Range("H34").Value = 2172.3286
$ws.Range("I34").Value = 1088.4615
$ws.Range("J34").Value = 3535.9033
$ws.Range("K34").Value = 1088.4615
$ws.Range("L34").Value = 3535.9033
$ws.Range("M34").Value = -886.4614999999999
$ws.Range("N34").Value = -3939.9033

$ws.Range("H132").Value = 905015.1
$ws.Range("I132").Value = 1288.4517
$ws.Range("K132").Value = 3865.3551
$ws.Range("M132").Value = -1335.3551

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1566.5555
$ws.Range("I35").Value = 550
$ws.Range("J35").Value = 1693.625
$ws.Range("K35").Value = 1650
$ws.Range("L35").Value = 5080.875
$ws.Range("M35").Value = -1362
$ws.Range("N35").Value = -5656.875

$ws.Range("H132").Value = 1089303.5
$ws.Range("I132").Value = 2292.9412
$ws.Range("J132").Value = 4169166.8
$ws.Range("K132").Value = 20636.4708
$ws.Range("L132").Value = 37522501.2
$ws.Range("M132").Value = -18106.4708
$ws.Range("N132").Value = -37527561.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 27779320
$ws.Range("I126").Value = 55556640
$ws.Range("K126").Value = 166669920
$ws.Range("M126").Value = -166667450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 956.2143
$ws.Range("I82").Value = 892.3333
$ws.Range("J82").Value = 973.63635
$ws.Range("K82").Value = 892.3333
$ws.Range("L82").Value = 973.63635
$ws.Range("M82").Value = -531.3333
$ws.Range("N82").Value = -1695.63635

$ws.Range("H85").Value = 956.2143
$ws.Range("I85").Value = 892.3333
$ws.Range("J85").Value = 973.63635
$ws.Range("K85").Value = 892.3333
$ws.Range("L85").Value = 973.63635
$ws.Range("M85").Value = 355.6667
$ws.Range("N85").Value = -3469.63635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1583.5
$ws.Range("I81").Value = 1775.25
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 3550.5
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -2489.5
$ws.Range("N81").Value = -4522

$ws.Range("H84").Value = 1583.5
$ws.Range("I84").Value = 1775.25
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 17752.5
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -12448.5
$ws.Range("N84").Value = -22608

Write-Host "Updated profit figures on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
